$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.124.51'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.552.26'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.33'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3821'
$ws.Range("E7").Value = '  +2.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3305'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.79'
$ws.Range("E9").Value = '  -9.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.134'
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07365'
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.13'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.815'
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("E17").Value = '  -3.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06639'
$ws.Range("E18").Value = '  -1.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.93'
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.360'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.05'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").Value = '22.120.84'
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.287'
$ws.Range("E25").Value = '  -4.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.537'
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.49'
$ws.Range("E27").Value = '  -1.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.13'
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.916'
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.087'
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.874'
$ws.Range("E33").Value = '  -4.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.896'
$ws.Range("E34").Value = '  -5.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08203'
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.299'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("E37").Value = '  -1.41%  '
$ws.Range("E40").Value = '  -5.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.236'
$ws.Range("E41").Value = '  -4.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.98'
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6048'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.80'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.731'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("E47").Value = '  -4.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.968'
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.86'
$ws.Range("E49").Value = '  -2.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.173'
$ws.Range("E50").Value = '  -3.32%  '
$ws.Range("E51").Value = '  -2.97%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.730'
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.558.42'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.25'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").Value = '1.740.93'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.301'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02320'
$ws.Range("E39").Value = '  -5.44%  '
